$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header C1: "Memory (GB)" -> "Heap Memory" ---
$ws.Range("C1").Value = "Heap Memory"

# --- Update C2 from numeric 29 to text "29 GB" ---
$ws.Range("C2").Value = "29 GB"

# --- Add new row 3 data (second memory-leak research observation, desktop run) ---
$ws.Range("A3").Value = 43401.775694444441
$ws.Range("A3").NumberFormat = "m/d/yy h:mm"

$ws.Range("B3").Value = 0.013865740740740739
$ws.Range("B3").NumberFormat = "h:mm:ss"

$ws.Range("C3").Value = "924.17196 GB"

$ws.Range("D3").Value = "Evaluated one individual on desktop with debug code single thread (break point at PushPG.compute_errors() line #38)."
$ws.Range("D3").WrapText = $true

# --- Row height for new row (auto height produced by wrapped text) ---
$ws.Rows.Item(3).RowHeight = 28.8

# --- Selection matches final authored state ---
[void]$ws.Range("D3").Select()
